# Carga_multimodal.xlsx update:
#  - E5 switches from the "ND" placeholder text to the actual figure (31452.1)
#  - The obsolete "ND No Disponible" legend row (row 34) is removed, shifting
#    the footnote rows below it up by one
#  - The "Actualización" footnote is bumped from Diciembre 2024 to Enero 2025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E5 now holds the real carretero figure instead of the "ND" text marker.
$ws.Range("E5").Value = 31452.1

# Delete the whole "ND No Disponible" row; everything beneath shifts up one row.
$ws.Rows(34).Delete()

# Update the "Actualización" note to the new month (now on row 33 after the shift).
$ws.Range("B33").Value = "Actualización: Enero 2025."
